$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.427210000000001
$ws.Range("H2").Value = 28.28163
$ws.Range("I2").Value = 0.2188083857550241
$ws.Range("J2").Value = 0.2188083857550241
$ws.Range("M2").Value = 0.4216986666666666
$ws.Range("N2").Value = 1.265096
$ws.Range("O2").Value = 0.2516921781598699
$ws.Range("P2").Value = 0.2516921781598699
$ws.Range("Q2").Value = 3.975441887386666
$ws.Range("R2").Value = 35.77897698648
$ws.Range("S2").Value = 0.05507235921032707
$ws.Range("T2").Value = 0.05507235921032707
$ws.Range("G3").Value = 9.427210000000001
$ws.Range("H3").Value = 28.28163
$ws.Range("I3").Value = 0.2188083857550241
$ws.Range("J3").Value = 0.2188083857550241
$ws.Range("O3").Value = 0.2613966125002536
$ws.Range("P3").Value = 0.2613966125002536
$ws.Range("Q3").Value = 4.12872203718
$ws.Range("R3").Value = 37.15849833462
$ws.Range("S3").Value = 0.05719577082301206
$ws.Range("T3").Value = 0.05719577082301206
$ws.Range("G4").Value = 9.427210000000001
$ws.Range("H4").Value = 28.28163
$ws.Range("I4").Value = 0.2188083857550241
$ws.Range("J4").Value = 0.2188083857550241
$ws.Range("M4").Value = 0.4328273333333333
$ws.Range("N4").Value = 1.298482
$ws.Range("O4").Value = 0.2583343579312433
$ws.Range("P4").Value = 0.2583343579312433
$ws.Range("Q4").Value = 4.080354165073333
$ws.Range("R4").Value = 36.72318748566
$ws.Range("S4").Value = 0.05652572384399595
$ws.Range("T4").Value = 0.05652572384399596
$ws.Range("G5").Value = 9.427210000000001
$ws.Range("H5").Value = 28.28163
$ws.Range("I5").Value = 0.2188083857550241
$ws.Range("J5").Value = 0.2188083857550241
$ws.Range("M5").Value = 0.38297
$ws.Range("N5").Value = 1.14891
$ws.Range("O5").Value = 0.2285768514086331
$ws.Range("P5").Value = 0.2285768514086331
$ws.Range("Q5").Value = 3.6103386137
$ws.Range("R5").Value = 32.4930475233
$ws.Range("S5").Value = 0.05001453187768902
$ws.Range("T5").Value = 0.05001453187768902
$ws.Range("I6").Value = 0.3808887290954196
$ws.Range("J6").Value = 0.3808887290954196
$ws.Range("M6").Value = 0.4216986666666666
$ws.Range("N6").Value = 1.265096
$ws.Range("O6").Value = 0.2516921781598699
$ws.Range("P6").Value = 0.2516921781598699
$ws.Range("Q6").Value = 6.920214702258665
$ws.Range("R6").Value = 62.28193232032799
$ws.Range("S6").Value = 0.09586671386257077
$ws.Range("T6").Value = 0.09586671386257077
$ws.Range("I7").Value = 0.3808887290954196
$ws.Range("J7").Value = 0.3808887290954196
$ws.Range("O7").Value = 0.2613966125002536
$ws.Range("P7").Value = 0.2613966125002536
$ws.Range("S7").Value = 0.09956302352506949
$ws.Range("T7").Value = 0.09956302352506949
$ws.Range("I8").Value = 0.3808887290954196
$ws.Range("J8").Value = 0.3808887290954196
$ws.Range("M8").Value = 0.4328273333333333
$ws.Range("N8").Value = 1.298482
$ws.Range("O8").Value = 0.2583343579312433
$ws.Range("P8").Value = 0.2583343579312433
$ws.Range("Q8").Value = 7.102839805847332
$ws.Range("R8").Value = 63.92555825262599
$ws.Range("S8").Value = 0.09839664527411247
$ws.Range("T8").Value = 0.0983966452741125
$ws.Range("I9").Value = 0.3808887290954196
$ws.Range("J9").Value = 0.3808887290954196
$ws.Range("M9").Value = 0.38297
$ws.Range("N9").Value = 1.14891
$ws.Range("O9").Value = 0.2285768514086331
$ws.Range("P9").Value = 0.2285768514086331
$ws.Range("Q9").Value = 6.284664463069999
$ws.Range("R9").Value = 56.56198016762999
$ws.Range("S9").Value = 0.08706234643366684
$ws.Range("T9").Value = 0.08706234643366684
$ws.Range("G10").Value = 7.213061
$ws.Range("H10").Value = 21.639183
$ws.Range("I10").Value = 0.1674173200514808
$ws.Range("J10").Value = 0.1674173200514808
$ws.Range("M10").Value = 0.4216986666666666
$ws.Range("N10").Value = 1.265096
$ws.Range("O10").Value = 0.2516921781598699
$ws.Range("P10").Value = 0.2516921781598699
$ws.Range("Q10").Value = 3.041738206285333
$ws.Range("R10").Value = 27.37564385656799
$ws.Range("S10").Value = 0.04213762994544525
$ws.Range("T10").Value = 0.04213762994544525
$ws.Range("G11").Value = 7.213061
$ws.Range("H11").Value = 21.639183
$ws.Range("I11").Value = 0.1674173200514808
$ws.Range("J11").Value = 0.1674173200514808
$ws.Range("O11").Value = 0.2613966125002536
$ws.Range("P11").Value = 0.2613966125002536
$ws.Range("Q11").Value = 3.159017769438
$ws.Range("R11").Value = 28.431159924942
$ws.Range("S11").Value = 0.04376232033532786
$ws.Range("T11").Value = 0.04376232033532786
$ws.Range("G12").Value = 7.213061
$ws.Range("H12").Value = 21.639183
$ws.Range("I12").Value = 0.1674173200514808
$ws.Range("J12").Value = 0.1674173200514808
$ws.Range("M12").Value = 0.4328273333333333
$ws.Range("N12").Value = 1.298482
$ws.Range("O12").Value = 0.2583343579312433
$ws.Range("P12").Value = 0.2583343579312433
$ws.Range("Q12").Value = 3.122009957800666
$ws.Range("R12").Value = 28.098089620206
$ws.Range("S12").Value = 0.04324964588206874
$ws.Range("T12").Value = 0.04324964588206875
$ws.Range("G13").Value = 7.213061
$ws.Range("H13").Value = 21.639183
$ws.Range("I13").Value = 0.1674173200514808
$ws.Range("J13").Value = 0.1674173200514808
$ws.Range("M13").Value = 0.38297
$ws.Range("N13").Value = 1.14891
$ws.Range("O13").Value = 0.2285768514086331
$ws.Range("P13").Value = 0.2285768514086331
$ws.Range("Q13").Value = 2.76238597117
$ws.Range("R13").Value = 24.86147374052999
$ws.Range("S13").Value = 0.03826772388863889
$ws.Range("T13").Value = 0.03826772388863889
$ws.Range("G14").Value = 10.03371566666667
$ws.Range("H14").Value = 30.101147
$ws.Range("I14").Value = 0.2328855650980756
$ws.Range("J14").Value = 0.2328855650980756
$ws.Range("M14").Value = 0.4216986666666666
$ws.Range("N14").Value = 1.265096
$ws.Range("O14").Value = 0.2516921781598699
$ws.Range("P14").Value = 0.2516921781598699
$ws.Range("Q14").Value = 4.231204518345778
$ws.Range("R14").Value = 38.080840665112
$ws.Range("S14").Value = 0.05861547514152682
$ws.Range("T14").Value = 0.05861547514152682
$ws.Range("G15").Value = 10.03371566666667
$ws.Range("H15").Value = 30.101147
$ws.Range("I15").Value = 0.2328855650980756
$ws.Range("J15").Value = 0.2328855650980756
$ws.Range("O15").Value = 0.2613966125002536
$ws.Range("P15").Value = 0.2613966125002536
$ws.Range("Q15").Value = 4.394346045942001
$ws.Range("R15").Value = 39.54911441347801
$ws.Range("S15").Value = 0.06087549781684426
$ws.Range("T15").Value = 0.06087549781684426
$ws.Range("G16").Value = 10.03371566666667
$ws.Range("H16").Value = 30.101147
$ws.Range("I16").Value = 0.2328855650980756
$ws.Range("J16").Value = 0.2328855650980756
$ws.Range("M16").Value = 0.4328273333333333
$ws.Range("N16").Value = 1.298482
$ws.Range("O16").Value = 0.2583343579312433
$ws.Range("P16").Value = 0.2583343579312433
$ws.Range("Q16").Value = 4.342866395428223
$ws.Range("R16").Value = 39.085797558854
$ws.Range("S16").Value = 0.06016234293106611
$ws.Range("T16").Value = 0.06016234293106613
$ws.Range("G17").Value = 10.03371566666667
$ws.Range("H17").Value = 30.101147
$ws.Range("I17").Value = 0.2328855650980756
$ws.Range("J17").Value = 0.2328855650980756
$ws.Range("M17").Value = 0.38297
$ws.Range("N17").Value = 1.14891
$ws.Range("O17").Value = 0.2285768514086331
$ws.Range("P17").Value = 0.2285768514086331
$ws.Range("Q17").Value = 3.842612088863333
$ws.Range("R17").Value = 34.58350879977
$ws.Range("S17").Value = 0.05323224920863837
$ws.Range("T17").Value = 0.05323224920863837
